$wb = $excel.ActiveWorkbook

$reb = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item('Assists'))
$reb.Name = 'Rebounds'

$rebHeader = New-Object 'object[,]' 1,16
$rebHeader[0,0] = 'Game Time (PST)'
$rebHeader[0,1] = 'Opponent'
$rebHeader[0,2] = 'Anthony Black'
$rebHeader[0,3] = 'Jonathan Isaac'
$rebHeader[0,4] = 'Tyus Jones'
$rebHeader[0,5] = 'Desmond Bane'
$rebHeader[0,6] = 'Jalen Suggs'
$rebHeader[0,7] = 'Paolo Banchero'
$rebHeader[0,8] = 'Orlando Robinson'
$rebHeader[0,9] = 'Jase Richardson'
$rebHeader[0,10] = 'Jett Howard'
$rebHeader[0,11] = 'Franz Wagner'
$rebHeader[0,12] = 'Tristan da Silva'
$rebHeader[0,13] = 'Wendell Carter Jr.'
$rebHeader[0,14] = 'Goga Bitadze'
$rebHeader[0,15] = 'Noah Penda'
$reb.Range('A1:P1').Value = $rebHeader
$reb.Range('A1:P1').Font.Bold = $true
$reb.Range('A1:P1').HorizontalAlignment = -4108
$reb.Range('A1:P1').VerticalAlignment = -4160
$reb.Range('A1:P1').Borders.LineStyle = 1

$reb.Range('A2:A12').NumberFormat = '@'

$rebData = New-Object 'object[,]' 11,16
$rebData[0,0] = '2025-10-22'
$rebData[0,1] = 'MIA'
$rebData[0,2] = 3
$rebData[0,3] = 1
$rebData[0,4] = 2
$rebData[0,5] = 4
$rebData[0,6] = 4
$rebData[0,7] = 11
$rebData[0,8] = 0
$rebData[0,9] = 0
$rebData[0,10] = 0
$rebData[0,11] = 4
$rebData[0,12] = 1
$rebData[0,13] = 8
$rebData[0,14] = 8
$rebData[0,15] = 0
$rebData[1,0] = '2025-10-24'
$rebData[1,1] = 'ATL'
$rebData[1,2] = 3
$rebData[1,3] = 0
$rebData[1,4] = 1
$rebData[1,5] = 7
$rebData[1,6] = 3
$rebData[1,7] = 6
$rebData[1,8] = 0
$rebData[1,9] = 0
$rebData[1,10] = 0
$rebData[1,11] = 6
$rebData[1,12] = 6
$rebData[1,13] = 8
$rebData[1,14] = 5
$rebData[1,15] = 0
$rebData[2,0] = '2025-10-25'
$rebData[2,1] = 'CHI'
$rebData[2,2] = 5
$rebData[2,3] = 1
$rebData[2,4] = 0
$rebData[2,5] = 5
$rebData[2,6] = 0
$rebData[2,7] = 10
$rebData[2,8] = 0
$rebData[2,9] = 0
$rebData[2,10] = 0
$rebData[2,11] = 7
$rebData[2,12] = 3
$rebData[2,13] = 13
$rebData[2,14] = 9
$rebData[2,15] = 0
$rebData[3,0] = '2025-10-27'
$rebData[3,1] = 'PHI'
$rebData[3,2] = 3
$rebData[3,3] = 0
$rebData[3,4] = 3
$rebData[3,5] = 3
$rebData[3,6] = 6
$rebData[3,7] = 7
$rebData[3,8] = 0
$rebData[3,9] = 0
$rebData[3,10] = 0
$rebData[3,11] = 4
$rebData[3,12] = 2
$rebData[3,13] = 10
$rebData[3,14] = 3
$rebData[3,15] = 2
$rebData[4,0] = '2025-10-29'
$rebData[4,1] = 'DET'
$rebData[4,2] = 4
$rebData[4,3] = 0
$rebData[4,4] = 0
$rebData[4,5] = 5
$rebData[4,6] = 2
$rebData[4,7] = 11
$rebData[4,8] = 0
$rebData[4,9] = 0
$rebData[4,10] = 1
$rebData[4,11] = 7
$rebData[4,12] = 3
$rebData[4,13] = 2
$rebData[4,14] = 6
$rebData[4,15] = 2
$rebData[5,0] = '2025-10-30'
$rebData[5,1] = 'CHA'
$rebData[5,2] = 1
$rebData[5,3] = 3
$rebData[5,4] = 0
$rebData[5,5] = 4
$rebData[5,6] = 0
$rebData[5,7] = 9
$rebData[5,8] = 0
$rebData[5,9] = 0
$rebData[5,10] = 0
$rebData[5,11] = 5
$rebData[5,12] = 3
$rebData[5,13] = 8
$rebData[5,14] = 3
$rebData[5,15] = 2
$rebData[6,0] = '2025-11-01'
$rebData[6,1] = 'WAS'
$rebData[6,2] = 2
$rebData[6,3] = 6
$rebData[6,4] = 2
$rebData[6,5] = 0
$rebData[6,6] = 2
$rebData[6,7] = 11
$rebData[6,8] = 0
$rebData[6,9] = 1
$rebData[6,10] = 1
$rebData[6,11] = 6
$rebData[6,12] = 2
$rebData[6,13] = 12
$rebData[6,14] = 6
$rebData[6,15] = 2
$rebData[7,0] = '2025-11-04'
$rebData[7,1] = 'ATL'
$rebData[7,2] = 4
$rebData[7,3] = 3
$rebData[7,4] = 0
$rebData[7,5] = 3
$rebData[7,6] = 2
$rebData[7,7] = 11
$rebData[7,8] = 0
$rebData[7,9] = 0
$rebData[7,10] = 0
$rebData[7,11] = 5
$rebData[7,12] = 4
$rebData[7,13] = 5
$rebData[7,14] = 2
$rebData[7,15] = 2
$rebData[8,0] = '2025-11-07'
$rebData[8,1] = 'BOS'
$rebData[8,2] = 4
$rebData[8,3] = 1
$rebData[8,4] = 1
$rebData[8,5] = 6
$rebData[8,6] = 3
$rebData[8,7] = 9
$rebData[8,8] = 0
$rebData[8,9] = 0
$rebData[8,10] = 0
$rebData[8,11] = 6
$rebData[8,12] = 2
$rebData[8,13] = 5
$rebData[8,14] = 7
$rebData[8,15] = 0
$rebData[9,0] = '2025-11-09'
$rebData[9,1] = 'BOS'
$rebData[9,2] = 3
$rebData[9,3] = 0
$rebData[9,4] = 0
$rebData[9,5] = 5
$rebData[9,6] = 8
$rebData[9,7] = 6
$rebData[9,8] = 0
$rebData[9,9] = 0
$rebData[9,10] = 0
$rebData[9,11] = 9
$rebData[9,12] = 5
$rebData[9,13] = 4
$rebData[9,14] = 2
$rebData[9,15] = 0
$rebData[10,0] = '2025-11-10'
$rebData[10,1] = 'POR'
$rebData[10,2] = 4
$rebData[10,3] = 6
$rebData[10,4] = 0
$rebData[10,5] = 3
$rebData[10,6] = 0
$rebData[10,7] = 9
$rebData[10,8] = 0
$rebData[10,9] = 0
$rebData[10,10] = 2
$rebData[10,11] = 9
$rebData[10,12] = 3
$rebData[10,13] = 9
$rebData[10,14] = 5
$rebData[10,15] = 0
$reb.Range('A2:P12').Value = $rebData

$pm3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item('Rebounds'))
$pm3.Name = '3PM'

$pm3Header = New-Object 'object[,]' 1,16
$pm3Header[0,0] = 'Game Time (PST)'
$pm3Header[0,1] = 'Opponent'
$pm3Header[0,2] = 'Anthony Black'
$pm3Header[0,3] = 'Jonathan Isaac'
$pm3Header[0,4] = 'Tyus Jones'
$pm3Header[0,5] = 'Desmond Bane'
$pm3Header[0,6] = 'Jalen Suggs'
$pm3Header[0,7] = 'Paolo Banchero'
$pm3Header[0,8] = 'Orlando Robinson'
$pm3Header[0,9] = 'Jase Richardson'
$pm3Header[0,10] = 'Jett Howard'
$pm3Header[0,11] = 'Franz Wagner'
$pm3Header[0,12] = 'Tristan da Silva'
$pm3Header[0,13] = 'Wendell Carter Jr.'
$pm3Header[0,14] = 'Goga Bitadze'
$pm3Header[0,15] = 'Noah Penda'
$pm3.Range('A1:P1').Value = $pm3Header
$pm3.Range('A1:P1').Font.Bold = $true
$pm3.Range('A1:P1').HorizontalAlignment = -4108
$pm3.Range('A1:P1').VerticalAlignment = -4160
$pm3.Range('A1:P1').Borders.LineStyle = 1

$pm3.Range('A2:A12').NumberFormat = '@'

$pm3Data = New-Object 'object[,]' 11,16
$pm3Data[0,0] = '2025-10-22'
$pm3Data[0,1] = 'MIA'
$pm3Data[0,2] = 1
$pm3Data[0,3] = 0
$pm3Data[0,4] = 0
$pm3Data[0,5] = 3
$pm3Data[0,6] = 2
$pm3Data[0,7] = 1
$pm3Data[0,8] = 0
$pm3Data[0,9] = 0
$pm3Data[0,10] = 0
$pm3Data[0,11] = 2
$pm3Data[0,12] = 3
$pm3Data[0,13] = 0
$pm3Data[0,14] = 0
$pm3Data[0,15] = 0
$pm3Data[1,0] = '2025-10-24'
$pm3Data[1,1] = 'ATL'
$pm3Data[1,2] = 0
$pm3Data[1,3] = 0
$pm3Data[1,4] = 0
$pm3Data[1,5] = 1
$pm3Data[1,6] = 3
$pm3Data[1,7] = 0
$pm3Data[1,8] = 0
$pm3Data[1,9] = 0
$pm3Data[1,10] = 1
$pm3Data[1,11] = 2
$pm3Data[1,12] = 2
$pm3Data[1,13] = 0
$pm3Data[1,14] = 0
$pm3Data[1,15] = 0
$pm3Data[2,0] = '2025-10-25'
$pm3Data[2,1] = 'CHI'
$pm3Data[2,2] = 1
$pm3Data[2,3] = 0
$pm3Data[2,4] = 0
$pm3Data[2,5] = 0
$pm3Data[2,6] = 0
$pm3Data[2,7] = 1
$pm3Data[2,8] = 0
$pm3Data[2,9] = 0
$pm3Data[2,10] = 0
$pm3Data[2,11] = 0
$pm3Data[2,12] = 0
$pm3Data[2,13] = 1
$pm3Data[2,14] = 0
$pm3Data[2,15] = 0
$pm3Data[3,0] = '2025-10-27'
$pm3Data[3,1] = 'PHI'
$pm3Data[3,2] = 1
$pm3Data[3,3] = 0
$pm3Data[3,4] = 0
$pm3Data[3,5] = 2
$pm3Data[3,6] = 1
$pm3Data[3,7] = 1
$pm3Data[3,8] = 0
$pm3Data[3,9] = 0
$pm3Data[3,10] = 0
$pm3Data[3,11] = 1
$pm3Data[3,12] = 0
$pm3Data[3,13] = 1
$pm3Data[3,14] = 0
$pm3Data[3,15] = 1
$pm3Data[4,0] = '2025-10-29'
$pm3Data[4,1] = 'DET'
$pm3Data[4,2] = 1
$pm3Data[4,3] = 0
$pm3Data[4,4] = 1
$pm3Data[4,5] = 1
$pm3Data[4,6] = 1
$pm3Data[4,7] = 0
$pm3Data[4,8] = 0
$pm3Data[4,9] = 0
$pm3Data[4,10] = 0
$pm3Data[4,11] = 3
$pm3Data[4,12] = 3
$pm3Data[4,13] = 2
$pm3Data[4,14] = 0
$pm3Data[4,15] = 0
$pm3Data[5,0] = '2025-10-30'
$pm3Data[5,1] = 'CHA'
$pm3Data[5,2] = 3
$pm3Data[5,3] = 0
$pm3Data[5,4] = 1
$pm3Data[5,5] = 0
$pm3Data[5,6] = 0
$pm3Data[5,7] = 2
$pm3Data[5,8] = 0
$pm3Data[5,9] = 0
$pm3Data[5,10] = 0
$pm3Data[5,11] = 3
$pm3Data[5,12] = 3
$pm3Data[5,13] = 3
$pm3Data[5,14] = 1
$pm3Data[5,15] = 0
$pm3Data[6,0] = '2025-11-01'
$pm3Data[6,1] = 'WAS'
$pm3Data[6,2] = 0
$pm3Data[6,3] = 1
$pm3Data[6,4] = 0
$pm3Data[6,5] = 1
$pm3Data[6,6] = 0
$pm3Data[6,7] = 3
$pm3Data[6,8] = 0
$pm3Data[6,9] = 0
$pm3Data[6,10] = 1
$pm3Data[6,11] = 2
$pm3Data[6,12] = 0
$pm3Data[6,13] = 2
$pm3Data[6,14] = 1
$pm3Data[6,15] = 1
$pm3Data[7,0] = '2025-11-04'
$pm3Data[7,1] = 'ATL'
$pm3Data[7,2] = 1
$pm3Data[7,3] = 0
$pm3Data[7,4] = 0
$pm3Data[7,5] = 2
$pm3Data[7,6] = 2
$pm3Data[7,7] = 0
$pm3Data[7,8] = 0
$pm3Data[7,9] = 0
$pm3Data[7,10] = 1
$pm3Data[7,11] = 0
$pm3Data[7,12] = 4
$pm3Data[7,13] = 1
$pm3Data[7,14] = 0
$pm3Data[7,15] = 0
$pm3Data[8,0] = '2025-11-07'
$pm3Data[8,1] = 'BOS'
$pm3Data[8,2] = 2
$pm3Data[8,3] = 1
$pm3Data[8,4] = 0
$pm3Data[8,5] = 2
$pm3Data[8,6] = 4
$pm3Data[8,7] = 1
$pm3Data[8,8] = 0
$pm3Data[8,9] = 0
$pm3Data[8,10] = 0
$pm3Data[8,11] = 2
$pm3Data[8,12] = 3
$pm3Data[8,13] = 2
$pm3Data[8,14] = 0
$pm3Data[8,15] = 0
$pm3Data[9,0] = '2025-11-09'
$pm3Data[9,1] = 'BOS'
$pm3Data[9,2] = 0
$pm3Data[9,3] = 0
$pm3Data[9,4] = 0
$pm3Data[9,5] = 0
$pm3Data[9,6] = 2
$pm3Data[9,7] = 0
$pm3Data[9,8] = 0
$pm3Data[9,9] = 0
$pm3Data[9,10] = 0
$pm3Data[9,11] = 1
$pm3Data[9,12] = 2
$pm3Data[9,13] = 2
$pm3Data[9,14] = 0
$pm3Data[9,15] = 0
$pm3Data[10,0] = '2025-11-10'
$pm3Data[10,1] = 'POR'
$pm3Data[10,2] = 1
$pm3Data[10,3] = 1
$pm3Data[10,4] = 0
$pm3Data[10,5] = 1
$pm3Data[10,6] = 0
$pm3Data[10,7] = 2
$pm3Data[10,8] = 0
$pm3Data[10,9] = 0
$pm3Data[10,10] = 1
$pm3Data[10,11] = 1
$pm3Data[10,12] = 0
$pm3Data[10,13] = 1
$pm3Data[10,14] = 0
$pm3Data[10,15] = 0
$pm3.Range('A2:P12').Value = $pm3Data

$avgReb = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item('Avg Assists'))
$avgReb.Name = 'Avg Rebounds'

$avgRebHeader = New-Object 'object[,]' 1,2
$avgRebHeader[0,0] = 'Player'
$avgRebHeader[0,1] = 'Avg Rebounds'
$avgReb.Range('A1:B1').Value = $avgRebHeader
$avgReb.Range('A1:B1').Font.Bold = $true
$avgReb.Range('A1:B1').HorizontalAlignment = -4108
$avgReb.Range('A1:B1').VerticalAlignment = -4160
$avgReb.Range('A1:B1').Borders.LineStyle = 1

$avgRebData = New-Object 'object[,]' 14,2
$avgRebData[0,0] = 'Paolo Banchero'
$avgRebData[0,1] = 9.090909090909092
$avgRebData[1,0] = 'Wendell Carter Jr.'
$avgRebData[1,1] = 7.636363636363637
$avgRebData[2,0] = 'Franz Wagner'
$avgRebData[2,1] = 6.181818181818182
$avgRebData[3,0] = 'Goga Bitadze'
$avgRebData[3,1] = 5.090909090909091
$avgRebData[4,0] = 'Desmond Bane'
$avgRebData[4,1] = 4.090909090909091
$avgRebData[5,0] = 'Jalen Suggs'
$avgRebData[5,1] = 3.75
$avgRebData[6,0] = 'Anthony Black'
$avgRebData[6,1] = 3.272727272727273
$avgRebData[7,0] = 'Tristan da Silva'
$avgRebData[7,1] = 3.090909090909091
$avgRebData[8,0] = 'Jonathan Isaac'
$avgRebData[8,1] = 2.1
$avgRebData[9,0] = 'Noah Penda'
$avgRebData[9,1] = 1.666666666666667
$avgRebData[10,0] = 'Tyus Jones'
$avgRebData[10,1] = 0.8181818181818182
$avgRebData[11,0] = 'Jett Howard'
$avgRebData[11,1] = 0.5714285714285714
$avgRebData[12,0] = 'Jase Richardson'
$avgRebData[12,1] = 0.1666666666666667
$avgRebData[13,0] = 'Orlando Robinson'
$avgRebData[13,1] = 0
$avgReb.Range('A2:B15').Value = $avgRebData

$avg3pm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item('Avg Rebounds'))
$avg3pm.Name = 'Avg 3PM'

$avg3pmHeader = New-Object 'object[,]' 1,2
$avg3pmHeader[0,0] = 'Player'
$avg3pmHeader[0,1] = 'Avg 3PM'
$avg3pm.Range('A1:B1').Value = $avg3pmHeader
$avg3pm.Range('A1:B1').Font.Bold = $true
$avg3pm.Range('A1:B1').HorizontalAlignment = -4108
$avg3pm.Range('A1:B1').VerticalAlignment = -4160
$avg3pm.Range('A1:B1').Borders.LineStyle = 1

$avg3pmData = New-Object 'object[,]' 14,2
$avg3pmData[0,0] = 'Jalen Suggs'
$avg3pmData[0,1] = 1.875
$avg3pmData[1,0] = 'Tristan da Silva'
$avg3pmData[1,1] = 1.818181818181818
$avg3pmData[2,0] = 'Franz Wagner'
$avg3pmData[2,1] = 1.545454545454545
$avg3pmData[3,0] = 'Wendell Carter Jr.'
$avg3pmData[3,1] = 1.363636363636364
$avg3pmData[4,0] = 'Desmond Bane'
$avg3pmData[4,1] = 1.181818181818182
$avg3pmData[5,0] = 'Anthony Black'
$avg3pmData[5,1] = 1
$avg3pmData[6,0] = 'Paolo Banchero'
$avg3pmData[6,1] = 1
$avg3pmData[7,0] = 'Jett Howard'
$avg3pmData[7,1] = 0.5714285714285714
$avg3pmData[8,0] = 'Noah Penda'
$avg3pmData[8,1] = 0.3333333333333333
$avg3pmData[9,0] = 'Jonathan Isaac'
$avg3pmData[9,1] = 0.3
$avg3pmData[10,0] = 'Tyus Jones'
$avg3pmData[10,1] = 0.1818181818181818
$avg3pmData[11,0] = 'Goga Bitadze'
$avg3pmData[11,1] = 0.1818181818181818
$avg3pmData[12,0] = 'Orlando Robinson'
$avg3pmData[12,1] = 0
$avg3pmData[13,0] = 'Jase Richardson'
$avg3pmData[13,1] = 0
$avg3pm.Range('A2:B15').Value = $avg3pmData

$wb.Worksheets.Item('Points').Activate()
